$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reshuffles the weekly price records across rows 2-14 (each row
# keeps the same columns A,B,C,E,F,G,N,O,Q,R which are constant anyway, but
# the per-record fields Fecha/Variedad/Calidad/Volumen/Precios move to a
# different row). Capture the "before" values per row, then write the
# "after" values computed from the target mapping.

$rows = 2..14

$fecha    = @{}
$variedad = @{}
$calidad  = @{}
$volumen  = @{}
$pmin     = @{}
$pmax     = @{}
$pprom    = @{}
$pkg      = @{}

foreach ($r in $rows) {
    $fecha[$r]    = $ws.Cells.Item($r, 4).Value2   # D - Fecha
    $variedad[$r] = $ws.Cells.Item($r, 8).Value2   # H - Variedad
    $calidad[$r]  = $ws.Cells.Item($r, 9).Value2   # I - Calidad
    $volumen[$r]  = $ws.Cells.Item($r, 10).Value2  # J - Volumen
    $pmin[$r]     = $ws.Cells.Item($r, 11).Value2  # K - Precio minimo
    $pmax[$r]     = $ws.Cells.Item($r, 12).Value2  # L - Precio maximo
    $pprom[$r]    = $ws.Cells.Item($r, 13).Value2  # M - Precio promedio ponderado
    $pkg[$r]      = $ws.Cells.Item($r, 16).Value2  # P - Precio $/Kg
}

# Target row <- source row (which row's original record now lands there)
$mapping = @{
    2  = 7
    3  = 13
    4  = 5
    5  = 6
    6  = 11
    7  = 4
    8  = 10
    9  = 14
    10 = 12
    11 = 2
    12 = 8
    13 = 3
    14 = 9
}

foreach ($target in $rows) {
    $src = $mapping[$target]
    $ws.Cells.Item($target, 4).Value  = $fecha[$src]
    $ws.Cells.Item($target, 8).Value  = $variedad[$src]
    $ws.Cells.Item($target, 9).Value  = $calidad[$src]
    $ws.Cells.Item($target, 10).Value = $volumen[$src]
    $ws.Cells.Item($target, 11).Value = $pmin[$src]
    $ws.Cells.Item($target, 12).Value = $pmax[$src]
    $ws.Cells.Item($target, 13).Value = $pprom[$src]
    $ws.Cells.Item($target, 16).Value = $pkg[$src]
}
